$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-coerced to a number by Excel
# (plain decimal look) need NumberFormat forced to Text first, to preserve
# the original text-cell representation (mirrors the OOXML inlineStr cells).
$textCells = @("D5","D6","D11","D14","D20","D24","D25","D26","D28","D32","D33","D36","D40","D41","D42","D47","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.783.10"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "3.107.24"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "585.97"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "145.46"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.102.04"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +7.55%  "
$ws.Range("D11").Value = "5.66"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "36.97"
$ws.Range("E14").Value = "  +4.09%  "
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "3.621.08"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "63.571.24"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "3.107.77"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("D20").Value = "462.38"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "12.99"
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("D25").Value = "81.26"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "9.23"
$ws.Range("E28").Value = "  +10.19%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "6.94"
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").Value = "26.74"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").Value = "0.0₃0859"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("D36").Value = "3.44"
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "50.42"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").Value = "438.37"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").Value = "8.69"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "2.882.66"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D47").Value = "36.57"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("D48").Value = "125.71"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "24.18"
$ws.Range("E51").Value = "  -1.31%  "
